# Weekly Fruta/Hortaliza update:
# Insert two new rows (a new week's worth of data, "Primera" and "Segunda"
# quality records) at the top of the data block that starts at row 200,
# pushing the existing rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 200 (shifts 200..215 -> 202..217)
$ws.Range("A200:A201").EntireRow.Insert()

# New row 200 - "Primera" quality for the new week (2022-01-17)
$ws.Cells.Item(200, 1).Value  = 9
$ws.Cells.Item(200, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(200, 3).Value  = "Metropolitana"
$ws.Cells.Item(200, 4).Value  = 44578
$ws.Cells.Item(200, 5).Value  = 13
$ws.Cells.Item(200, 6).Value  = 100112043
$ws.Cells.Item(200, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(200, 8).Value  = "Sin especificar"
$ws.Cells.Item(200, 9).Value  = "Primera"
$ws.Cells.Item(200, 10).Value = 79
$ws.Cells.Item(200, 11).Value = 9000
$ws.Cells.Item(200, 12).Value = 10000
$ws.Cells.Item(200, 13).Value = 9494
$ws.Cells.Item(200, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(200, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(200, 16).Value = 158
$ws.Cells.Item(200, 17).Value = 60
$ws.Cells.Item(200, 18).Value = "Hortaliza"

# New row 201 - "Segunda" quality for the new week (2022-01-17)
$ws.Cells.Item(201, 1).Value  = 9
$ws.Cells.Item(201, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(201, 3).Value  = "Metropolitana"
$ws.Cells.Item(201, 4).Value  = 44578
$ws.Cells.Item(201, 5).Value  = 13
$ws.Cells.Item(201, 6).Value  = 100112043
$ws.Cells.Item(201, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(201, 8).Value  = "Sin especificar"
$ws.Cells.Item(201, 9).Value  = "Segunda"
$ws.Cells.Item(201, 10).Value = 43
$ws.Cells.Item(201, 11).Value = 7000
$ws.Cells.Item(201, 12).Value = 7000
$ws.Cells.Item(201, 13).Value = 7000
$ws.Cells.Item(201, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(201, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(201, 16).Value = 70
$ws.Cells.Item(201, 17).Value = 100
$ws.Cells.Item(201, 18).Value = "Hortaliza"
